$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether | Ether
$ws.Cells.Item(15, 8).Value = 2632.7083
$ws.Cells.Item(15, 9).Value = 2632.7083
$ws.Cells.Item(15, 11).Value = 7898.124899999999
$ws.Cells.Item(15, 13).Value = -7729.124899999999

# Row 43: Growing Is Knowing | Growth Formula Gamma
$ws.Cells.Item(43, 8).Value = 2909.0588
$ws.Cells.Item(43, 9).Value = 3207.3333
$ws.Cells.Item(43, 10).Value = 2746.3635
$ws.Cells.Item(43, 11).Value = 3207.3333
$ws.Cells.Item(43, 12).Value = 2746.3635
$ws.Cells.Item(43, 13).Value = -3138.3333
$ws.Cells.Item(43, 14).Value = -2884.3635

# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Cells.Item(76, 8).Value = 2500
$ws.Cells.Item(76, 9).Value = 2500
$ws.Cells.Item(76, 11).Value = 2500
$ws.Cells.Item(76, 13).Value = -2185

# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Cells.Item(79, 8).Value = 2500
$ws.Cells.Item(79, 9).Value = 2500
$ws.Cells.Item(79, 11).Value = 2500
$ws.Cells.Item(79, 13).Value = -1408

# Row 103: Let Loose the Juice | Persimmon Tannin
$ws.Cells.Item(103, 8).Value = 2760.1
$ws.Cells.Item(103, 9).Value = 3184.1667
$ws.Cells.Item(103, 10).Value = 2124
$ws.Cells.Item(103, 11).Value = 9552.500100000001
$ws.Cells.Item(103, 12).Value = 6372
$ws.Cells.Item(103, 13).Value = -8966.500100000001
$ws.Cells.Item(103, 14).Value = -7544

# Row 127: Liquid Competence | Competent Craftsman's Draught
$ws.Cells.Item(127, 8).Value = 1402.1666
$ws.Cells.Item(127, 9).Value = 1273.8
$ws.Cells.Item(127, 11).Value = 3821.4
$ws.Cells.Item(127, 13).Value = 1138.6

# Row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Cells.Item(135, 8).Value = 641.9091
$ws.Cells.Item(135, 9).Value = 478.8
$ws.Cells.Item(135, 10).Value = 777.8333
$ws.Cells.Item(135, 11).Value = 4309.2
$ws.Cells.Item(135, 12).Value = 7000.4997
$ws.Cells.Item(135, 13).Value = -1774.2
$ws.Cells.Item(135, 14).Value = -12070.4997

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Cells.Item(138, 8).Value = 267201.5
$ws.Cells.Item(138, 9).Value = 4259.2905
$ws.Cells.Item(138, 10).Value = 394564.12
$ws.Cells.Item(138, 11).Value = 12777.8715
$ws.Cells.Item(138, 12).Value = 1183692.36
$ws.Cells.Item(138, 13).Value = -7637.871500000001
$ws.Cells.Item(138, 14).Value = -1193972.36

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Cells.Item(141, 8).Value = 4845
$ws.Cells.Item(141, 9).Value = 4739.6665
$ws.Cells.Item(141, 10).Value = 5266.3335
$ws.Cells.Item(141, 11).Value = 14218.9995
$ws.Cells.Item(141, 12).Value = 15799.0005
$ws.Cells.Item(141, 13).Value = -9038.999500000002
$ws.Cells.Item(141, 14).Value = -26159.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Cells.Item(2, 8).Value = 883.5833
$ws.Cells.Item(2, 9).Value = 894.7778
$ws.Cells.Item(2, 10).Value = 850
$ws.Cells.Item(2, 11).Value = 894.7778
$ws.Cells.Item(2, 12).Value = 850
$ws.Cells.Item(2, 13).Value = -781.7778
$ws.Cells.Item(2, 14).Value = -1076

# Row 32: Ingot We Trust | Steel Ingot
$ws.Cells.Item(32, 8).Value = 4360.7085
$ws.Cells.Item(32, 9).Value = 3872.0725
$ws.Cells.Item(32, 11).Value = 3872.0725
$ws.Cells.Item(32, 13).Value = -3585.0725

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Cells.Item(61, 8).Value = 7168.0586
$ws.Cells.Item(61, 9).Value = 2042.75
$ws.Cells.Item(61, 11).Value = 2042.75
$ws.Cells.Item(61, 13).Value = -1830.75

# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Cells.Item(102, 8).Value = 4968.56
$ws.Cells.Item(102, 9).Value = 5004.9473
$ws.Cells.Item(102, 10).Value = 4853.3335
$ws.Cells.Item(102, 11).Value = 5004.9473
$ws.Cells.Item(102, 12).Value = 4853.3335
$ws.Cells.Item(102, 13).Value = -3382.9473
$ws.Cells.Item(102, 14).Value = -8097.3335

# Row 116: No Scope | Titanbronze Ingot
$ws.Cells.Item(116, 8).Value = 883.5833
$ws.Cells.Item(116, 9).Value = 894.7778
$ws.Cells.Item(116, 10).Value = 850
$ws.Cells.Item(116, 11).Value = 894.7778
$ws.Cells.Item(116, 12).Value = 850
$ws.Cells.Item(116, 13).Value = 1399.2222
$ws.Cells.Item(116, 14).Value = -5438

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 1964.7778
$ws.Cells.Item(132, 9).Value = 1439.6444
$ws.Cells.Item(132, 11).Value = 4318.933199999999
$ws.Cells.Item(132, 13).Value = -1788.933199999999

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value = 7168.0586
$ws.Cells.Item(136, 9).Value = 2042.75
$ws.Cells.Item(136, 11).Value = 6128.25
$ws.Cells.Item(136, 13).Value = -3578.25

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Cells.Item(3, 8).Value = 883.5833
$ws.Cells.Item(3, 9).Value = 894.7778
$ws.Cells.Item(3, 10).Value = 850
$ws.Cells.Item(3, 11).Value = 894.7778
$ws.Cells.Item(3, 12).Value = 850
$ws.Cells.Item(3, 13).Value = -780.7778
$ws.Cells.Item(3, 14).Value = -1078

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Cells.Item(99, 8).Value = 3614.2
$ws.Cells.Item(99, 9).Value = 1010
$ws.Cells.Item(99, 10).Value = 4265.25
$ws.Cells.Item(99, 11).Value = 1010
$ws.Cells.Item(99, 12).Value = 4265.25
$ws.Cells.Item(99, 13).Value = 488
$ws.Cells.Item(99, 14).Value = -7261.25

# Row 141: Awl Dreams Come True | Ra'Kaznar Awl
$ws.Cells.Item(141, 8).Value = 74994.164
$ws.Cells.Item(141, 10).Value = 79993
$ws.Cells.Item(141, 12).Value = 79993
$ws.Cells.Item(141, 14).Value = -90353

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Cells.Item(31, 8).Value = 5152.186
$ws.Cells.Item(31, 9).Value = 4237.2
$ws.Cells.Item(31, 10).Value = 5947.826
$ws.Cells.Item(31, 11).Value = 4237.2
$ws.Cells.Item(31, 12).Value = 5947.826
$ws.Cells.Item(31, 13).Value = -3942.2
$ws.Cells.Item(31, 14).Value = -6537.826

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Cells.Item(34, 8).Value = 5152.186
$ws.Cells.Item(34, 9).Value = 4237.2
$ws.Cells.Item(34, 10).Value = 5947.826
$ws.Cells.Item(34, 11).Value = 4237.2
$ws.Cells.Item(34, 12).Value = 5947.826
$ws.Cells.Item(34, 13).Value = -4035.2
$ws.Cells.Item(34, 14).Value = -6351.826

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Cells.Item(58, 8).Value = 3996.4375
$ws.Cells.Item(58, 9).Value = 2744.125
$ws.Cells.Item(58, 11).Value = 2744.125
$ws.Cells.Item(58, 13).Value = -2541.125

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Cells.Item(132, 8).Value = 3677.658
$ws.Cells.Item(132, 9).Value = 2980.5925
$ws.Cells.Item(132, 10).Value = 5388.636
$ws.Cells.Item(132, 11).Value = 8941.7775
$ws.Cells.Item(132, 12).Value = 16165.908
$ws.Cells.Item(132, 13).Value = -6411.7775
$ws.Cells.Item(132, 14).Value = -21225.908

# Row 135: The Wing's Wings | Ceiba Wings
$ws.Cells.Item(135, 8).Value = 120665.555
$ws.Cells.Item(135, 10).Value = 120665.555
$ws.Cells.Item(135, 12).Value = 120665.555
$ws.Cells.Item(135, 14).Value = -130805.555

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Cells.Item(136, 8).Value = 3996.4375
$ws.Cells.Item(136, 9).Value = 2744.125
$ws.Cells.Item(136, 11).Value = 8232.375
$ws.Cells.Item(136, 13).Value = -5682.375

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water | Boiled Egg
$ws.Cells.Item(4, 8).Value = 6273478.5
$ws.Cells.Item(4, 9).Value = 5900684
$ws.Cells.Item(4, 10).Value = 8044251
$ws.Cells.Item(4, 11).Value = 17702052
$ws.Cells.Item(4, 12).Value = 24132753
$ws.Cells.Item(4, 13).Value = -17701940
$ws.Cells.Item(4, 14).Value = -24132977

# Row 14: Keep Your Powder Dry | Kukuru Powder
$ws.Cells.Item(14, 8).Value = 495.5
$ws.Cells.Item(14, 9).Value = 495.5
$ws.Cells.Item(14, 11).Value = 1486.5
$ws.Cells.Item(14, 13).Value = -1313.5

# Row 56: Culture Club | Crowned Pie
$ws.Cells.Item(56, 8).Value = 8061.25
$ws.Cells.Item(56, 9).Value = 8061.25
$ws.Cells.Item(56, 11).Value = 8061.25
$ws.Cells.Item(56, 13).Value = -7531.25

# Row 122: Salt of the North | Northern Sea Salt
$ws.Cells.Item(122, 8).Value = 1138.2693
$ws.Cells.Item(122, 9).Value = 1067.2
$ws.Cells.Item(122, 10).Value = 1155.1904
$ws.Cells.Item(122, 11).Value = 9604.800000000001
$ws.Cells.Item(122, 12).Value = 10396.7136
$ws.Cells.Item(122, 13).Value = -7154.800000000001
$ws.Cells.Item(122, 14).Value = -15296.7136

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Cells.Item(80, 8).Value = 55557732
$ws.Cells.Item(80, 9).Value = 90911010
$ws.Cells.Item(80, 10).Value = 2592.2856
$ws.Cells.Item(80, 11).Value = 90911010
$ws.Cells.Item(80, 12).Value = 2592.2856
$ws.Cells.Item(80, 13).Value = -90910012
$ws.Cells.Item(80, 14).Value = -4588.2856

# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Cells.Item(83, 8).Value = 55557732
$ws.Cells.Item(83, 9).Value = 90911010
$ws.Cells.Item(83, 10).Value = 2592.2856
$ws.Cells.Item(83, 11).Value = 454555050
$ws.Cells.Item(83, 12).Value = 12961.428
$ws.Cells.Item(83, 13).Value = -454550058
$ws.Cells.Item(83, 14).Value = -22945.428

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Cells.Item(102, 8).Value = 1484.5
$ws.Cells.Item(102, 9).Value = 1292.1111
$ws.Cells.Item(102, 11).Value = 1292.1111
$ws.Cells.Item(102, 13).Value = 329.8888999999999

# Row 132: On Board for Lar | Lar Ingot
$ws.Cells.Item(132, 8).Value = 3365
$ws.Cells.Item(132, 9).Value = 2695.5789
$ws.Cells.Item(132, 10).Value = 4636.9
$ws.Cells.Item(132, 11).Value = 8086.736699999999
$ws.Cells.Item(132, 12).Value = 13910.7
$ws.Cells.Item(132, 13).Value = -5556.736699999999
$ws.Cells.Item(132, 14).Value = -18970.7

$ws = $wb.Worksheets.Item("LTW")
# Row 9: From the Sands to the Stage | Leather Himantes
$ws.Cells.Item(9, 8).Value = 5600
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 13).ClearContents()

# Row 16: Saddle Sore | Hard Leather
$ws.Cells.Item(16, 8).Value = 756.1579
$ws.Cells.Item(16, 9).Value = 685.64703
$ws.Cells.Item(16, 11).Value = 685.64703
$ws.Cells.Item(16, 13).Value = -515.64703

# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Cells.Item(22, 8).Value = 590.5
$ws.Cells.Item(22, 9).Value = 587.375
$ws.Cells.Item(22, 10).Value = 594.6667
$ws.Cells.Item(22, 11).Value = 587.375
$ws.Cells.Item(22, 12).Value = 594.6667
$ws.Cells.Item(22, 13).Value = -292.375
$ws.Cells.Item(22, 14).Value = -1184.6667

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Cells.Item(27, 8).Value = 590.5
$ws.Cells.Item(27, 9).Value = 587.375
$ws.Cells.Item(27, 10).Value = 594.6667
$ws.Cells.Item(27, 11).Value = 587.375
$ws.Cells.Item(27, 12).Value = 594.6667
$ws.Cells.Item(27, 13).Value = -480.375
$ws.Cells.Item(27, 14).Value = -808.6667

# Row 40: Best Served Toad | Toad Leather
$ws.Cells.Item(40, 8).Value = 4199.354
$ws.Cells.Item(40, 9).Value = 4155.3555
$ws.Cells.Item(40, 11).Value = 4155.3555
$ws.Cells.Item(40, 13).Value = -4019.3555

# Row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws.Cells.Item(55, 8).Value = 367
$ws.Cells.Item(55, 9).Value = 314.14285
$ws.Cells.Item(55, 11).Value = 314.14285
$ws.Cells.Item(55, 13).Value = -141.14285

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Cells.Item(136, 8).Value = 8690.615
$ws.Cells.Item(136, 9).Value = 7854.143
$ws.Cells.Item(136, 11).Value = 23562.429
$ws.Cells.Item(136, 13).Value = -21012.429

# Row 138: Freezing Toes | Gomphotherium Boots of Striking
$ws.Cells.Item(138, 8).Value = 69407
$ws.Cells.Item(138, 10).Value = 69407
$ws.Cells.Item(138, 12).Value = 69407
$ws.Cells.Item(138, 14).Value = -79687

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value = 2578.78
$ws.Cells.Item(132, 9).Value = 2443.068
$ws.Cells.Item(132, 11).Value = 7329.204000000001
$ws.Cells.Item(132, 13).Value = -4799.204000000001

# Row 138: Halfgloves, Full Effort | Rroneek Serge Halfgloves of Healing
$ws.Cells.Item(138, 8).Value = 99844.5
$ws.Cells.Item(138, 10).Value = 99844.5
$ws.Cells.Item(138, 12).Value = 99844.5
$ws.Cells.Item(138, 14).Value = -110124.5
